$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 3.55
$ws.Range("J2").Value = 3.15
$ws.Range("T2").Value = 1.86
$ws.Range("U2").Value = 1.95
$ws.Range("F3").Value = 1.9
$ws.Range("H3").Value = 4.4
$ws.Range("I3").Value = 4.9
$ws.Range("L3").Value = 1.36
$ws.Range("T3").Value = 1.87
$ws.Range("F4").Value = 2.7
$ws.Range("G4").Value = 3.85
$ws.Range("H4").Value = 2.46
$ws.Range("I4").Value = 3.45
$ws.Range("J4").Value = 2.68
$ws.Range("K4").Value = 4.5
$ws.Range("L4").Value = 1.36
$ws.Range("P4").Value = 1.57
$ws.Range("V4").Value = 1.41
$ws.Range("W4").Value = 1.35
$ws.Range("F5").Value = 1.54
$ws.Range("G5").Value = 1.71
$ws.Range("J5").Value = 3.6
$ws.Range("K5").Value = 4.7
$ws.Range("M5").Value = 1.06
$ws.Range("R5").Value = 1.36
$ws.Range("T5").Value = 1.93
$ws.Range("W5").Value = 2.4
$ws.Range("AB5").Value = 9.8
$ws.Range("F6").Value = 2.94
$ws.Range("G6").Value = 3.2
$ws.Range("H6").Value = 2.42
$ws.Range("I6").Value = 2.66
$ws.Range("J6").Value = 3.3
$ws.Range("N6").Value = 3.6
$ws.Range("P6").Value = 1.87
$ws.Range("Q6").Value = 1.92
$ws.Range("R6").Value = 1.34
$ws.Range("S6").Value = 3.3
$ws.Range("T6").Value = 1.72
$ws.Range("V6").Value = 1.6
$ws.Range("W6").Value = 1.46
$ws.Range("X6").Value = 17.5
$ws.Range("Z6").Value = 17.5
$ws.Range("AA6").Value = 980
$ws.Range("AB6").Value = 13
$ws.Range("AE6").Value = 980
$ws.Range("AF6").Value = 23
$ws.Range("AG6").Value = 14
$ws.Range("AO6").Value = 27
$ws.Range("G7").Value = 1.48
$ws.Range("T7").Value = 1.9
$ws.Range("AB7").Value = 9.6
$ws.Range("F8").Value = 1.57
$ws.Range("G8").Value = 1.71
$ws.Range("H8").Value = 6.2
$ws.Range("I8").Value = 10.5
$ws.Range("J8").Value = 3.2
$ws.Range("K8").Value = 4.4
$ws.Range("L8").Value = 1.44
$ws.Range("N8").Value = 2.44
$ws.Range("O8").Value = 1.47
$ws.Range("P8").Value = 1.57
$ws.Range("Q8").Value = 2.2
$ws.Range("R8").Value = 1.2
$ws.Range("S8").Value = 1.05
$ws.Range("T8").Value = 2.28
$ws.Range("U8").Value = 1.61
$ws.Range("W8").Value = 2.4
$ws.Range("AB8").Value = 990
$ws.Range("F9").Value = 1.91
$ws.Range("G9").Value = 2.3
$ws.Range("H9").Value = 3.4
$ws.Range("I9").Value = 5.1
$ws.Range("J9").Value = 2.7
$ws.Range("K9").Value = 5
$ws.Range("L9").Value = 1.26
$ws.Range("N9").Value = 2.6
$ws.Range("P9").Value = 1.87
$ws.Range("Q9").Value = 1.62
$ws.Range("R9").Value = 1.38
$ws.Range("S9").Value = 2.28
$ws.Range("V9").Value = 1.24
$ws.Range("W9").Value = 1.77
